# Auto-generated edit script: updates Asura_Profits crafting-profit cached values
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 671.41174  # H28
$ws.Cells.Item(28, 9).Value = 697.75  # I28
$ws.Cells.Item(28, 10).Value = 250  # J28
$ws.Cells.Item(28, 11).Value = 697.75  # K28
$ws.Cells.Item(28, 12).Value = 250  # L28
$ws.Cells.Item(28, 13).Value = -212.75  # M28
$ws.Cells.Item(28, 14).Value = -1220  # N28
$ws.Cells.Item(61, 8).Value = 443.8  # H61
$ws.Cells.Item(61, 9).Value = 443.8  # I61
$ws.Cells.Item(61, 10).Value = 0  # J61
$ws.Cells.Item(61, 11).Value = 1331.4  # K61
$ws.Cells.Item(61, 12).Value = 0  # L61
$ws.Cells.Item(61, 13).Value = -1159.4  # M61
$ws.Cells.Item(61, 14).Value = ""  # N61: clear (was -1595)
$ws.Cells.Item(76, 8).Value = 4642.857  # H76
$ws.Cells.Item(76, 9).Value = 5400  # I76
$ws.Cells.Item(76, 10).Value = 4340  # J76
$ws.Cells.Item(76, 11).Value = 5400  # K76
$ws.Cells.Item(76, 12).Value = 4340  # L76
$ws.Cells.Item(76, 13).Value = -5085  # M76
$ws.Cells.Item(76, 14).Value = -4970  # N76
$ws.Cells.Item(79, 8).Value = 4642.857  # H79
$ws.Cells.Item(79, 9).Value = 5400  # I79
$ws.Cells.Item(79, 10).Value = 4340  # J79
$ws.Cells.Item(79, 11).Value = 5400  # K79
$ws.Cells.Item(79, 12).Value = 4340  # L79
$ws.Cells.Item(79, 13).Value = -4308  # M79
$ws.Cells.Item(79, 14).Value = -6524  # N79
$ws.Cells.Item(98, 8).Value = 6577.115  # H98
$ws.Cells.Item(98, 9).Value = 5114.476  # I98
$ws.Cells.Item(98, 10).Value = 12720.2  # J98
$ws.Cells.Item(98, 11).Value = 5114.476  # K98
$ws.Cells.Item(98, 12).Value = 12720.2  # L98
$ws.Cells.Item(98, 13).Value = -3616.476  # M98
$ws.Cells.Item(98, 14).Value = -15716.2  # N98
$ws.Cells.Item(122, 8).Value = 6577.115  # H122
$ws.Cells.Item(122, 9).Value = 5114.476  # I122
$ws.Cells.Item(122, 10).Value = 12720.2  # J122
$ws.Cells.Item(122, 11).Value = 15343.428  # K122
$ws.Cells.Item(122, 12).Value = 38160.60000000001  # L122
$ws.Cells.Item(122, 13).Value = -12893.428  # M122
$ws.Cells.Item(122, 14).Value = -43060.60000000001  # N122
$ws.Cells.Item(129, 8).Value = 953.4545000000001  # H129
$ws.Cells.Item(129, 10).Value = 1038.5217  # J129
$ws.Cells.Item(129, 12).Value = 3115.5651  # L129
$ws.Cells.Item(129, 14).Value = -13115.5651  # N129
$ws.Cells.Item(132, 8).Value = 1590.0416  # H132
$ws.Cells.Item(132, 9).Value = 1288.3125  # I132
$ws.Cells.Item(132, 11).Value = 3864.9375  # K132
$ws.Cells.Item(132, 13).Value = -1334.9375  # M132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10443.964  # H32
$ws.Cells.Item(32, 9).Value = 10911.632  # I32
$ws.Cells.Item(32, 10).Value = 6001.125  # J32
$ws.Cells.Item(32, 11).Value = 10911.632  # K32
$ws.Cells.Item(32, 12).Value = 6001.125  # L32
$ws.Cells.Item(32, 13).Value = -10624.632  # M32
$ws.Cells.Item(32, 14).Value = -6575.125  # N32
$ws.Cells.Item(107, 8).Value = 33000  # H107
$ws.Cells.Item(107, 10).Value = 33000  # J107
$ws.Cells.Item(107, 12).Value = 33000  # L107
$ws.Cells.Item(107, 14).Value = -40680  # N107
$ws.Cells.Item(122, 8).Value = 1922.091  # H122
$ws.Cells.Item(122, 9).Value = 1889.2069  # I122
$ws.Cells.Item(122, 10).Value = 2160.5  # J122
$ws.Cells.Item(122, 11).Value = 5667.620699999999  # K122
$ws.Cells.Item(122, 12).Value = 6481.5  # L122
$ws.Cells.Item(122, 13).Value = -3217.620699999999  # M122
$ws.Cells.Item(122, 14).Value = -11381.5  # N122
$ws.Cells.Item(131, 8).Value = 24500  # H131
$ws.Cells.Item(131, 10).Value = 24500  # J131
$ws.Cells.Item(131, 12).Value = 24500  # L131
$ws.Cells.Item(131, 14).Value = -34580  # N131
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 24612.5  # H81
$ws.Cells.Item(81, 10).Value = 24612.5  # J81
$ws.Cells.Item(81, 12).Value = 24612.5  # L81
$ws.Cells.Item(81, 14).Value = -26734.5  # N81
$ws.Cells.Item(84, 8).Value = 24612.5  # H84
$ws.Cells.Item(84, 10).Value = 24612.5  # J84
$ws.Cells.Item(84, 12).Value = 73837.5  # L84
$ws.Cells.Item(84, 14).Value = -84445.5  # N84
$ws.Cells.Item(94, 8).Value = 1163.5625  # H94
$ws.Cells.Item(94, 9).Value = 985.6087  # I94
$ws.Cells.Item(94, 11).Value = 985.6087  # K94
$ws.Cells.Item(94, 13).Value = -534.6087  # M94
$ws.Cells.Item(109, 8).Value = 34059.668  # H109
$ws.Cells.Item(109, 10).Value = 34059.668  # J109
$ws.Cells.Item(109, 12).Value = 34059.668  # L109
$ws.Cells.Item(109, 14).Value = -36833.668  # N109
$ws.Cells.Item(132, 8).Value = 75939.8  # H132
$ws.Cells.Item(132, 10).Value = 75939.8  # J132
$ws.Cells.Item(132, 12).Value = 75939.8  # L132
$ws.Cells.Item(132, 14).Value = -86059.8  # N132
$ws.Cells.Item(133, 8).Value = 66391.664  # H133
$ws.Cells.Item(133, 10).Value = 66391.664  # J133
$ws.Cells.Item(133, 12).Value = 66391.664  # L133
$ws.Cells.Item(133, 14).Value = -76511.664  # N133
$ws.Cells.Item(134, 8).Value = 335583.66  # H134
$ws.Cells.Item(134, 9).Value = 456204.66  # I134
$ws.Cells.Item(134, 10).Value = 3875.875  # J134
$ws.Cells.Item(134, 11).Value = 1368613.98  # K134
$ws.Cells.Item(134, 12).Value = 11627.625  # L134
$ws.Cells.Item(134, 13).Value = -1366078.98  # M134
$ws.Cells.Item(134, 14).Value = -16697.625  # N134
$ws.Cells.Item(135, 8).Value = 60735  # H135
$ws.Cells.Item(135, 10).Value = 60735  # J135
$ws.Cells.Item(135, 12).Value = 60735  # L135
$ws.Cells.Item(135, 14).Value = -70875  # N135
$ws.Cells.Item(137, 8).Value = 0  # H137
$ws.Cells.Item(137, 10).Value = 0  # J137
$ws.Cells.Item(137, 12).Value = 0  # L137
$ws.Cells.Item(137, 14).Value = ""  # N137: clear (was -60008.332)
$ws.Cells.Item(138, 8).Value = 39644.832  # H138
$ws.Cells.Item(138, 10).Value = 39644.832  # J138
$ws.Cells.Item(138, 12).Value = 39644.832  # L138
$ws.Cells.Item(138, 14).Value = -49924.832  # N138
$ws.Cells.Item(139, 8).Value = 179940  # H139
$ws.Cells.Item(139, 10).Value = 179940  # J139
$ws.Cells.Item(139, 12).Value = 179940  # L139
$ws.Cells.Item(139, 14).Value = -190220  # N139
$ws.Cells.Item(140, 8).Value = 89593.336  # H140
$ws.Cells.Item(140, 10).Value = 89593.336  # J140
$ws.Cells.Item(140, 12).Value = 89593.336  # L140
$ws.Cells.Item(140, 14).Value = -99953.336  # N140
$ws.Cells.Item(141, 8).Value = 44000  # H141
$ws.Cells.Item(141, 10).Value = 44000  # J141
$ws.Cells.Item(141, 12).Value = 44000  # L141
$ws.Cells.Item(141, 14).Value = -54360  # N141
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 274643  # H28
$ws.Cells.Item(28, 10).Value = 274643  # J28
$ws.Cells.Item(28, 12).Value = 274643  # L28
$ws.Cells.Item(28, 14).Value = -275133  # N28
$ws.Cells.Item(31, 8).Value = 11907007  # H31
$ws.Cells.Item(31, 9).Value = 18869426  # I31
$ws.Cells.Item(31, 10).Value = 3518.742  # J31
$ws.Cells.Item(31, 11).Value = 18869426  # K31
$ws.Cells.Item(31, 12).Value = 3518.742  # L31
$ws.Cells.Item(31, 13).Value = -18869131  # M31
$ws.Cells.Item(31, 14).Value = -4108.742  # N31
$ws.Cells.Item(34, 8).Value = 11907007  # H34
$ws.Cells.Item(34, 9).Value = 18869426  # I34
$ws.Cells.Item(34, 10).Value = 3518.742  # J34
$ws.Cells.Item(34, 11).Value = 18869426  # K34
$ws.Cells.Item(34, 12).Value = 3518.742  # L34
$ws.Cells.Item(34, 13).Value = -18869224  # M34
$ws.Cells.Item(34, 14).Value = -3922.742  # N34
$ws.Cells.Item(58, 8).Value = 1174.1708  # H58
$ws.Cells.Item(58, 9).Value = 1150.2354  # I58
$ws.Cells.Item(58, 10).Value = 1290.4286  # J58
$ws.Cells.Item(58, 11).Value = 1150.2354  # K58
$ws.Cells.Item(58, 12).Value = 1290.4286  # L58
$ws.Cells.Item(58, 13).Value = -947.2354  # M58
$ws.Cells.Item(58, 14).Value = -1696.4286  # N58
$ws.Cells.Item(70, 8).Value = 11975  # H70
$ws.Cells.Item(70, 10).Value = 11975  # J70
$ws.Cells.Item(70, 12).Value = 11975  # L70
$ws.Cells.Item(70, 14).Value = -12605  # N70
$ws.Cells.Item(73, 8).Value = 11975  # H73
$ws.Cells.Item(73, 10).Value = 11975  # J73
$ws.Cells.Item(73, 12).Value = 11975  # L73
$ws.Cells.Item(73, 14).Value = -14159  # N73
$ws.Cells.Item(123, 8).Value = 65647.62  # H123
$ws.Cells.Item(123, 10).Value = 65647.62  # J123
$ws.Cells.Item(123, 12).Value = 65647.62  # L123
$ws.Cells.Item(123, 14).Value = -75447.62  # N123
$ws.Cells.Item(132, 8).Value = 2677.3333  # H132
$ws.Cells.Item(132, 9).Value = 2419.2964  # I132
$ws.Cells.Item(132, 10).Value = 4999.6665  # J132
$ws.Cells.Item(132, 11).Value = 7257.889200000001  # K132
$ws.Cells.Item(132, 12).Value = 14998.9995  # L132
$ws.Cells.Item(132, 13).Value = -4727.889200000001  # M132
$ws.Cells.Item(132, 14).Value = -20058.9995  # N132
$ws.Cells.Item(133, 8).Value = 67109.5  # H133
$ws.Cells.Item(133, 10).Value = 67109.5  # J133
$ws.Cells.Item(133, 12).Value = 67109.5  # L133
$ws.Cells.Item(133, 14).Value = -72169.5  # N133
$ws.Cells.Item(134, 8).Value = 1140.75  # H134
$ws.Cells.Item(134, 9).Value = 960.8857400000001  # I134
$ws.Cells.Item(134, 10).Value = 1840.2222  # J134
$ws.Cells.Item(134, 11).Value = 2882.65722  # K134
$ws.Cells.Item(134, 12).Value = 5520.6666  # L134
$ws.Cells.Item(134, 13).Value = -347.6572200000001  # M134
$ws.Cells.Item(134, 14).Value = -10590.6666  # N134
$ws.Cells.Item(136, 8).Value = 1174.1708  # H136
$ws.Cells.Item(136, 9).Value = 1150.2354  # I136
$ws.Cells.Item(136, 10).Value = 1290.4286  # J136
$ws.Cells.Item(136, 11).Value = 3450.7062  # K136
$ws.Cells.Item(136, 12).Value = 3871.2858  # L136
$ws.Cells.Item(136, 13).Value = -900.7062000000001  # M136
$ws.Cells.Item(136, 14).Value = -8971.2858  # N136
$ws.Cells.Item(138, 8).Value = 52331.934  # H138
$ws.Cells.Item(138, 10).Value = 52331.934  # J138
$ws.Cells.Item(138, 12).Value = 52331.934  # L138
$ws.Cells.Item(138, 14).Value = -62611.934  # N138
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 920231.5  # H12
$ws.Cells.Item(12, 9).Value = 143.66667  # I12
$ws.Cells.Item(12, 10).Value = 1073579.5  # J12
$ws.Cells.Item(12, 11).Value = 431.00001  # K12
$ws.Cells.Item(12, 12).Value = 3220738.5  # L12
$ws.Cells.Item(12, 13).Value = -258.00001  # M12
$ws.Cells.Item(12, 14).Value = -3221084.5  # N12
$ws.Cells.Item(113, 8).Value = 793.0625  # H113
$ws.Cells.Item(113, 10).Value = 853.38464  # J113
$ws.Cells.Item(113, 12).Value = 2560.15392  # L113
$ws.Cells.Item(113, 14).Value = -6900.15392  # N113
$ws.Cells.Item(117, 8).Value = 37087.645  # H117
$ws.Cells.Item(117, 10).Value = 39915.152  # J117
$ws.Cells.Item(117, 12).Value = 119745.456  # L117
$ws.Cells.Item(117, 14).Value = -126629.456  # N117
$ws.Cells.Item(131, 8).Value = 867.62  # H131
$ws.Cells.Item(131, 10).Value = 887.38947  # J131
$ws.Cells.Item(131, 12).Value = 2662.16841  # L131
$ws.Cells.Item(131, 14).Value = -12742.16841  # N131
$ws.Cells.Item(132, 8).Value = 1619.6666  # H132
$ws.Cells.Item(132, 9).Value = 800.1429000000001  # I132
$ws.Cells.Item(132, 10).Value = 2141.182  # J132
$ws.Cells.Item(132, 11).Value = 7201.2861  # K132
$ws.Cells.Item(132, 12).Value = 19270.638  # L132
$ws.Cells.Item(132, 13).Value = -4671.2861  # M132
$ws.Cells.Item(132, 14).Value = -24330.638  # N132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(87, 8).Value = 30000  # H87
$ws.Cells.Item(87, 10).Value = 30000  # J87
$ws.Cells.Item(87, 12).Value = 30000  # L87
$ws.Cells.Item(87, 14).Value = -32496  # N87
$ws.Cells.Item(90, 8).Value = 30000  # H90
$ws.Cells.Item(90, 10).Value = 30000  # J90
$ws.Cells.Item(90, 12).Value = 90000  # L90
$ws.Cells.Item(90, 14).Value = -102480  # N90
$ws.Cells.Item(109, 8).Value = 21035  # H109
$ws.Cells.Item(109, 10).Value = 21035  # J109
$ws.Cells.Item(109, 12).Value = 21035  # L109
$ws.Cells.Item(109, 14).Value = -23115  # N109
$ws.Cells.Item(122, 8).Value = 2922.0977  # H122
$ws.Cells.Item(122, 9).Value = 3017.8823  # I122
$ws.Cells.Item(122, 11).Value = 9053.6469  # K122
$ws.Cells.Item(122, 13).Value = -6603.6469  # M122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4047.8333  # H40
$ws.Cells.Item(40, 9).Value = 4047.8333  # I40
$ws.Cells.Item(40, 11).Value = 4047.8333  # K40
$ws.Cells.Item(40, 13).Value = -3911.8333  # M40
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 15627672  # H122
$ws.Cells.Item(122, 9).Value = 19232366  # I122
$ws.Cells.Item(122, 10).Value = 7328.3335  # J122
$ws.Cells.Item(122, 11).Value = 57697098  # K122
$ws.Cells.Item(122, 12).Value = 21985.0005  # L122
$ws.Cells.Item(122, 13).Value = -57694648  # M122
$ws.Cells.Item(122, 14).Value = -26885.0005  # N122
